# Generate Report for Handoff
# - Update the "Status"-style cells from "Handed back: in sync with en-US" to
#   "Ready for handoff"
# - Bump the related handoff-generation timestamps by one minute
# - Narrow the columns that used to hold the long "Handed back..." text down
#   to a width appropriate for the shorter "Ready for handoff" label

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Target column width (in characters) we want stored in the sheet once Excel
# rounds it to its internal pixel grid. The host snaps ColumnWidth to steps
# of 1/6 character, so we dial the input in slightly to land the stored
# value as close as possible to 17.2159881591797.
$targetColWidth = 16.3826548258464

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-12 19:14:21"

$wsOverview.Range("E:E").ColumnWidth = $targetColWidth
$wsOverview.Range("F:F").ColumnWidth = $targetColWidth

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-12 19:14:13"

$wsZhCn.Range("C:C").ColumnWidth = $targetColWidth

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-12 19:14:21"

$wsDeDe.Range("C:C").ColumnWidth = $targetColWidth
